$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/number, week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# --- Weekly crime-stat table updates (rows 15-29) ---
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 2
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -50
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("I15").Value = 1
$ws.Range("I15").NumberFormat = "#,##0"
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = -66.666666666666
$ws.Range("L15").Value = -66.666666666666
$ws.Range("L15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -75
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 47
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 161.111111111111
$ws.Range("I16").Value = 32
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = 128.571428571429
$ws.Range("L16").Value = 33.333333333333
$ws.Range("M16").Value = 45.454545454545
$ws.Range("N16").Value = -66.315789473684
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -64.285714285714
$ws.Range("F17").Value = 37
$ws.Range("H17").Value = -36.206896551724
$ws.Range("I17").Value = 27
$ws.Range("J17").Value = 42
$ws.Range("K17").Value = -35.714285714285
$ws.Range("L17").Value = 17.391304347826
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = -27.027027027027
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 500
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 177.777777777778
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = 6
$ws.Range("K18").Value = 300
$ws.Range("L18").Value = 60
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = -68.831168831168
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -23.076923076923
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = 16.666666666666
$ws.Range("I19").Value = 26
$ws.Range("J19").Value = 29
$ws.Range("K19").Value = -10.344827586206
$ws.Range("L19").Value = -10.344827586206
$ws.Range("M19").Value = 52.941176470588
$ws.Range("N19").Value = -16.129032258064
$ws.Range("C20").Value = "'0"
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = -33.333333333333
$ws.Range("J20").Value = 18
$ws.Range("K20").Value = -27.777777777777
$ws.Range("L20").Value = 8.333333333333
$ws.Range("M20").Value = 85.714285714285
$ws.Range("N20").Value = -69.767441860465
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -36.363636363636
$ws.Range("F21").Value = 172
$ws.Range("G21").Value = 157
$ws.Range("H21").Value = 9.554140127388
$ws.Range("I21").Value = 123
$ws.Range("J21").Value = 112
$ws.Range("K21").Value = 9.821428571428
$ws.Range("L21").Value = 14.953271028037
$ws.Range("M21").Value = 89.230769230769
$ws.Range("N21").Value = -57.586206896551
$ws.Range("D23").Value = "'0"
$ws.Range("E23").Value = "'***.*"
$ws.Range("G23").Value = 2
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -35.714285714285
$ws.Range("F24").Value = 72
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = -19.101123595505
$ws.Range("I24").Value = 51
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = -15
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 24.390243902439
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -7.142857142857
$ws.Range("F25").Value = 63
$ws.Range("H25").Value = -5.970149253731
$ws.Range("I25").Value = 51
$ws.Range("J25").Value = 42
$ws.Range("K25").Value = 21.428571428571
$ws.Range("L25").Value = 15.90909090909
$ws.Range("M25").Value = 24.390243902439
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 4
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = -75
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -83.333333333333
$ws.Range("I26").Value = 1
$ws.Range("I26").NumberFormat = "#,##0"
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = -80
$ws.Range("L26").Value = -75
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = -75
$ws.Range("L27").Value = -33.333333333333
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F28").Value = "'0"
$ws.Range("H28").Value = -100
$ws.Range("J28").Value = 1
$ws.Range("J28").NumberFormat = "#,##0"
$ws.Range("K28").Value = -100
$ws.Range("K28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F29").Value = "'0"
$ws.Range("H29").Value = -100
$ws.Range("J29").Value = 1
$ws.Range("J29").NumberFormat = "#,##0"
$ws.Range("K29").Value = -100
$ws.Range("K29").NumberFormat = "#,##0.0;""-""#,##0.0"
